$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the empty paragraph immediately following the code line that sets
# className="postWrapper" - this is where the new "Like functionality"
# commentary needs to be inserted (after the 2nd blank paragraph that
# follows the code snippet).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "postWrapper") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find anchor paragraph containing 'postWrapper'"
}

# The anchor paragraph is followed by two blank paragraphs, then the new
# "//Like functionality with useState hook" paragraph is inserted, then one
# blank paragraph, then the new explanatory "*React useState hook..."
# paragraph.
$firstBlank = $targetIndex + 2

# --- Insert paragraph 1: "//Like functionality with useState hook" ---
$p = $d.Paragraphs.Item($firstBlank)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item($firstBlank + 1)
$insertRange1 = $newPara1.Range
$insertRange1.Collapse(0)
$xml1 = "<w:p $ns><w:r><w:lastRenderedPageBreak/><w:tab/></w:r><w:r><w:tab/><w:t>//Like functionality with useState hook</w:t></w:r></w:p>"
$null = $insertRange1.InsertXML($xml1)

# --- Insert paragraph 2: "*React useState hook allows us..." ---
$secondBlank = $firstBlank + 2
$p2 = $d.Paragraphs.Item($secondBlank)
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Item($secondBlank + 1)
$insertRange2 = $newPara2.Range
$insertRange2.Collapse(0)
$xml2 = "<w:p $ns><w:r><w:t>*React useState hook allows us to track state in a function component. State generally refer to the data or properties that need to be trackin in an application.</w:t></w:r></w:p>"
$null = $insertRange2.InsertXML($xml2)
